$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 16: was an empty/placeholder template row -> becomes a real diary
# entry (1/30 class session).  Copy number/date formats from an existing
# data row (row 10) so the per-column styles (date/time/text/mood) match.
# ---------------------------------------------------------------------
$ws.Range("A10:G10").Copy() | Out-Null
$ws.Range("A16:G16").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(16).RowHeight = 29.85

$ws.Cells.Item(16,1).Value = 43860
$ws.Cells.Item(16,2).Value = 0.708333333333333
$ws.Cells.Item(16,3).Value = "Class"
$ws.Cells.Item(16,4).Value = "Attend lecture"
$ws.Cells.Item(16,5).Value = "Learned more about UML diagrams and other models"
$ws.Cells.Item(16,6).Value = "Choosing the right representation for the job is important"
$ws.Cells.Item(16,7).Value = "Positive"

# ---------------------------------------------------------------------
# Row 17: likewise becomes a real diary entry (2/1 team session).
# ---------------------------------------------------------------------
$ws.Range("A10:G10").Copy() | Out-Null
$ws.Range("A17:G17").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(17).RowHeight = 29.85

$ws.Cells.Item(17,1).Value = 43862
$ws.Cells.Item(17,2).Value = 0.541666666666667
$ws.Cells.Item(17,3).Value = "Team"
$ws.Cells.Item(17,4).Value = "Write packet for 2 essential features"
$ws.Cells.Item(17,5).Value = "All goals"
$ws.Cells.Item(17,6).Value = "One function often branches out and touches a lot of different other parts."
$ws.Cells.Item(17,7).Value = "Positive"

# ---------------------------------------------------------------------
# Rows 18-21 keep their existing placeholder-template formatting/values,
# they just shift position within the shared-string table (handled
# automatically).  Nothing else to change on them.
#
# Rows 22-25 were empty trailing rows; they become additional copies of
# the placeholder template row (same values/formats as row 18), with a
# slightly taller row height (15.75 instead of 15.5).
# ---------------------------------------------------------------------
$ws.Range("A18:G18").Copy() | Out-Null
foreach ($r in 22..25) {
    $target = $ws.Range("A" + $r + ":G" + $r)
    $target.PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r,1).Value = "<what day?>"
    $ws.Cells.Item($r,2).Value = "<what time?>"
    $ws.Cells.Item($r,3).Value = "<as applicable, with whom?>"
    $ws.Cells.Item($r,4).Value = "<what did you want to accomplish?>"
    $ws.Cells.Item($r,5).Value = "<what did you actually accomplish?>"
    $ws.Cells.Item($r,6).Value = "<what insight(s) did you gain?>"
    $ws.Cells.Item($r,7).Value = "<how did you feel during the activity?>"
    $ws.Rows.Item($r).RowHeight = 15.75
}

# ---------------------------------------------------------------------
# Sheet view bookkeeping: the used range grew from A1:G22 to A1:G26 and
# the visible/selected cell moved from G15 to G18.
# ---------------------------------------------------------------------
$ws.Cells.Item(18,7).Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 8
